# Apply cell-value updates from the crypto price refresh.
# Source cells are stored as text (inlineStr) in the workbook even when
# the text looks like a plain decimal number (e.g. "1.01"), so for those
# cells we force the Text number format first to stop Excel's COM layer
# from auto-converting the assigned string into a numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.019.63"
$ws.Range("D3").Value = "1.561.14"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.19"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.248"
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0854"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "1.784.24"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "1.562.40"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "26.996.89"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.09"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.37"
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  +1.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.19"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.93"
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.06"
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0474"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("E31").Value = "  +3.59%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  +3.35%  "
$ws.Range("D34").Value = "1.428.64"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("E36").Value = "  +9.44%  "
$ws.Range("E37").Value = "  +2.24%  "
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.532"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("E40").Value = "  +2.06%  "
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.64"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.74"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("D47").Value = "1.696.88"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.79"
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").Value = "  +3.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0959"
$ws.Range("E51").Value = "  +0.51%  "
